# Updated Links & Pdf's — Outer PDF's are added; CD Data-Flow-Analysis links were modified.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CD")

# --- New content, written in the order the original author typed it so that
# --- the shared-string table grows in the same sequence as the target file. ---

# Intermediate Code / Three Address Code block (rows 12-13)
$ws.Range("A12").Value = "Intermediate Code"
$ws.Range("B12").Value = "Three Address code"
$ws.Range("C12").Value = "https://www.youtube.com/watch?v=yFVCw0N0nxo&ab_channel=SudhakarAtchala"
$ws.Range("C13").Value = "https://www.youtube.com/watch?v=sv5Qmq8Hjz4&t=140s&ab_channel=SudhakarAtchala"
$ws.Range("B13").Value = "Quadrples, Triples,Indirect Triples"

# Code Optimization block (rows 15-17), with a merged, styled empty cell at C16
$ws.Range("A15").Value = "Code Optimization"
$ws.Range("C15").Value = "https://www.youtube.com/watch?v=yHZFVz6TVmI&list=PLXj4XH7LcRfC9pGMWuM6UWE3V4YZ9TZzM&index=51&ab_channel=SudhakarAtchala"
$ws.Range("B15").Value = "Constant Propagation"
$ws.Range("B16").Value = "Common Subexpression Elimination"
$ws.Range("C17").Value = "https://www.gatevidyalay.com/code-optimization-techniques/"
$ws.Range("B17").Value = "Reading Material"

# Data-Flow-Analysis block (row 19)
$ws.Range("B19").Value = "Data-Flow-Analysis"
$ws.Range("C19").Value = "https://www.youtube.com/watch?v=OROXJ9-wUQE&t=1600s&ab_channel=MayurNaik"

# Tutorials Point block (row 5) — reuses the existing "Youtube Playlist" string
$ws.Range("B5").Value = "Tutorials Point"
$ws.Range("C5").Value = "https://youtu.be/0hU5-aLtaxo"
$ws.Range("A5").Value = "Youtube Playlist"

# C15:C16 is a merged cell, left-aligned / vertically centered (style reused from sheet)
$ws.Range("C15").VerticalAlignment = -4108
$ws.Range("C15").HorizontalAlignment = -4131
$ws.Range("C16").VerticalAlignment = -4108
$ws.Range("C16").HorizontalAlignment = -4131
$ws.Range("C15:C16").Merge()

# Column B widened to fit the new, longer topic names
$ws.Range("B1").ColumnWidth = 29.25

# Final selection as left by the author
[void]$ws.Range("C5").Select()
